$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.890262829142614
$ws.Range("C2").Value = 0.01399785247282637
$ws.Range("D2").Value = 0.02396514811884032
$ws.Range("E2").Value = 0.06982531799436043
$ws.Range("F2").Value = 7.206549594585169
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("J2").Value = 0.2626523037579247
$ws.Range("K2").Value = 1.335197723484129
$ws.Range("L2").Value = 0.2571189953132063
# Row 3
$ws.Range("B3").Value = 1.885041995309194
$ws.Range("C3").Value = 0.0122322069236489
$ws.Range("D3").Value = 0.02095857312384908
$ws.Range("E3").Value = 0.07035925692298939
$ws.Range("F3").Value = 7.013907244349809
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("J3").Value = 0.2592143794182391
$ws.Range("K3").Value = 1.327649784914968
$ws.Range("L3").Value = 0.2592806340793103
# Row 4
$ws.Range("B4").Value = 1.883239106679895
$ws.Range("C4").Value = 0.01116520683416411
$ws.Range("D4").Value = 0.01911154582446528
$ws.Range("E4").Value = 0.07071429031197507
$ws.Range("F4").Value = 6.896144708780099
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("J4").Value = 0.2571212288252838
$ws.Range("K4").Value = 1.324156127388022
$ws.Range("L4").Value = 0.2607903568574983
# Row 5
$ws.Range("B5").Value = 1.88285706661901
$ws.Range("C5").Value = 0.01073459623033557
$ws.Range("D5").Value = 0.01835847611946662
$ws.Range("E5").Value = 0.07086582266493302
$ws.Range("F5").Value = 6.848282981670849
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("J5").Value = 0.2562726549881518
$ws.Range("K5").Value = 1.323019033808052
$ws.Range("L5").Value = 0.2614514928433849
# Row 6
$ws.Range("B6").Value = 1.882814926617357
$ws.Range("C6").Value = 0.01066334438681338
$ws.Range("D6").Value = 0.01823340160761688
$ws.Range("E6").Value = 0.07089139888144302
$ws.Range("F6").Value = 6.840343176378099
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("J6").Value = 0.2561320140509267
$ws.Range("K6").Value = 1.322847522790497
$ws.Range("L6").Value = 0.2615640477969308
# Row 7
$ws.Range("B7").Value = 1.883232526581764
$ws.Range("C7").Value = 0.01115938259624016
$ws.Range("D7").Value = 0.0191013914369762
$ws.Range("E7").Value = 0.07071630616068703
$ws.Range("F7").Value = 6.895498717293492
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("J7").Value = 0.257109766916173
$ws.Range("K7").Value = 1.324139632064174
$ws.Range("L7").Value = 0.2607990872357462
# Row 8
$ws.Range("B8").Value = 1.88817149144532
$ws.Range("C8").Value = 0.01338544342741699
$ws.Range("D8").Value = 0.02292859594171404
$ws.Range("E8").Value = 0.07000378557478903
$ws.Range("F8").Value = 7.140016304967162
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("J8").Value = 0.261463185694069
$ws.Range("K8").Value = 1.332358247934792
$ws.Range("L8").Value = 0.2578264735719387
# Row 9
$ws.Range("B9").Value = 1.908995095785912
$ws.Range("C9").Value = 0.01789140105790921
$ws.Range("D9").Value = 0.03043280216797939
$ws.Range("E9").Value = 0.0688216180199861
$ws.Range("F9").Value = 7.623826277718052
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("J9").Value = 0.2701442728344645
$ws.Range("K9").Value = 1.357543038645247
$ws.Range("L9").Value = 0.2534439857367516
# Row 10
$ws.Range("B10").Value = 1.931103118597207
$ws.Range("C10").Value = 0.02129476794834773
$ws.Range("D10").Value = 0.0359554411162577
$ws.Range("E10").Value = 0.0680832922879322
$ws.Range("F10").Value = 7.982195614809797
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("J10").Value = 0.2766152031089604
$ws.Range("K10").Value = 1.381603391434368
$ws.Range("L10").Value = 0.2511052198034562
# Row 11
$ws.Range("B11").Value = 1.942643926733496
$ws.Range("C11").Value = 0.02286477156737021
$ws.Range("D11").Value = 0.0384720049939915
$ws.Range("E11").Value = 0.06777549425172857
$ws.Range("F11").Value = 8.145924044504341
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("J11").Value = 0.2795803299336228
$ws.Range("K11").Value = 1.393762649370871
$ws.Range("L11").Value = 0.2502324226437054
# Row 12
$ws.Range("B12").Value = 1.947227798893749
$ws.Range("C12").Value = 0.02346254730001363
$ws.Range("D12").Value = 0.03942574452744907
$ws.Range("E12").Value = 0.0676629607096686
$ws.Range("F12").Value = 8.208029526195389
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("J12").Value = 0.280706313691816
$ws.Range("K12").Value = 1.398542112831052
$ws.Range("L12").Value = 0.2499293876078212
# Row 13
$ws.Range("B13").Value = 1.946231074973781
$ws.Range("C13").Value = 0.02333365922866903
$ws.Range("D13").Value = 0.03922030287152722
$ws.Range("E13").Value = 0.06768701811590283
$ws.Range("F13").Value = 8.194649265328223
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("J13").Value = 0.2804636715588558
$ws.Range("K13").Value = 1.397504978871495
$ws.Range("L13").Value = 0.2499934298102033
# Row 14
$ws.Range("B14").Value = 1.943016761884309
$ws.Range("C14").Value = 0.02291388514673542
$ws.Range("D14").Value = 0.03855045336908347
$ws.Range("E14").Value = 0.06776615550165133
$ws.Range("F14").Value = 8.151031372085242
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("J14").Value = 0.2796729016828152
$ws.Range("K14").Value = 1.394152348574437
$ws.Range("L14").Value = 0.2502069411843308
# Row 15
$ws.Range("B15").Value = 1.941075729684457
$ws.Range("C15").Value = 0.02265718802156869
$ws.Range("D15").Value = 0.03814025627691819
$ws.Range("E15").Value = 0.06781515293611129
$ws.Range("F15").Value = 8.124327950491988
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("J15").Value = 0.2791889451382659
$ws.Range("K15").Value = 1.392121571163017
$ws.Range("L15").Value = 0.2503413008148172
# Row 16
$ws.Range("B16").Value = 1.930378778080978
$ws.Range("C16").Value = 0.02119261291797159
$ws.Range("D16").Value = 0.03579107724745256
$ws.Range("E16").Value = 0.0681039712925946
$ws.Range("F16").Value = 7.971510074006062
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("J16").Value = 0.2764218628409765
$ws.Range("K16").Value = 1.380833222359371
$ws.Range("L16").Value = 0.2511661041947946
# Row 17
$ws.Range("B17").Value = 1.924196746640234
$ws.Range("C17").Value = 0.02029980736904236
$ws.Range("D17").Value = 0.03435114011524831
$ws.Range("E17").Value = 0.06828833171704751
$ws.Range("F17").Value = 7.877944353102578
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("J17").Value = 0.2747298937350635
$ws.Range("K17").Value = 1.37421944387421
$ws.Range("L17").Value = 0.251721037970924
# Row 18
$ws.Range("B18").Value = 1.920780647475311
$ws.Range("C18").Value = 0.01978833369192046
$ws.Range("D18").Value = 0.03352332063802521
$ws.Range("E18").Value = 0.06839701384165409
$ws.Range("F18").Value = 7.824193899484612
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("J18").Value = 0.2737587344547805
$ws.Range("K18").Value = 1.370529621670642
$ws.Range("L18").Value = 0.2520582106314322
# Row 19
$ws.Range("B19").Value = 1.919647991112356
$ws.Range("C19").Value = 0.01961550514023713
$ws.Range("D19").Value = 0.03324309819934967
$ws.Range("E19").Value = 0.06843426607478609
$ws.Range("F19").Value = 7.806006153242578
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("J19").Value = 0.273430260921117
$ws.Range("K19").Value = 1.369299920617294
$ws.Range("L19").Value = 0.2521754614885623
# Row 20
$ws.Range("B20").Value = 1.924840380287264
$ws.Range("C20").Value = 0.02039463561749244
$ws.Range("D20").Value = 0.03450438195275751
$ws.Range("E20").Value = 0.06826843280577233
$ws.Range("F20").Value = 7.8878977181582
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("J20").Value = 0.2749097975344839
$ws.Range("K20").Value = 1.374911664575677
$ws.Range("L20").Value = 0.251660102556535
# Row 21
$ws.Range("B21").Value = 1.943955083868815
$ws.Range("C21").Value = 0.02303709385596164
$ws.Range("D21").Value = 0.03874718223167406
$ws.Range("E21").Value = 0.06774280186001214
$ws.Range("F21").Value = 8.163840122444071
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("J21").Value = 0.2799050838069945
$ws.Range("K21").Value = 1.395132343927571
$ws.Range("L21").Value = 0.2501434821306248
# Row 22
$ws.Range("B22").Value = 1.957692913152755
$ws.Range("C22").Value = 0.02478308781913086
$ws.Range("D22").Value = 0.04152466629237495
$ws.Range("E22").Value = 0.06742271399872024
$ws.Range("F22").Value = 8.344798704672883
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("J22").Value = 0.2831882114252267
$ws.Range("K22").Value = 1.409368070028933
$ws.Range("L22").Value = 0.2493124217796066
# Row 23
$ws.Range("B23").Value = 1.950246743489004
$ws.Range("C23").Value = 0.02384944147628687
$ws.Range("D23").Value = 0.04004180135838453
$ws.Range("E23").Value = 0.06759141036617677
$ws.Range("F23").Value = 8.248160235388355
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("J23").Value = 0.2814342365340181
$ws.Range("K23").Value = 1.401676685860821
$ws.Range("L23").Value = 0.2497413238884221
# Row 24
$ws.Range("B24").Value = 1.924548963574807
$ws.Range("C24").Value = 0.02035175813009005
$ws.Range("D24").Value = 0.0344351012699633
$ws.Range("E24").Value = 0.06827742072100573
$ws.Range("F24").Value = 7.883397672114455
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("J24").Value = 0.2748284581238991
$ws.Range("K24").Value = 1.374598361118188
$ws.Range("L24").Value = 0.2516875949717345
# Row 25
$ws.Range("B25").Value = 1.90216746134206
$ws.Range("C25").Value = 0.01665663344046919
$ws.Range("D25").Value = 0.02840179777228968
$ws.Range("E25").Value = 0.06911849401364911
$ws.Range("F25").Value = 7.492451898675114
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("J25").Value = 0.2677798612286608
$ws.Range("K25").Value = 1.349756155385478
$ws.Range("L25").Value = 0.2544747810555634

Write-Host "done"
